{"js": "// Update the answer text in each populated table cell of the practice\n// worksheet (\"two-digit number divided by one-digit number\") in document\n// order, mapping each old \"a\u00f7b=c, d\" string to its replacement.\nconst pairs = [\n  [\"11\u00f77=1, 4\", \"35\u00f77=5, 0\"],\n  [\"14\u00f77=2, 0\", \"45\u00f73=15, 0\"],\n  [\"89\u00f79=9, 8\", \"67\u00f79=7, 4\"],\n  [\"60\u00f77=8, 4\", \"13\u00f79=1, 4\"],\n  [\"66\u00f74=16, 2\", \"70\u00f79=7, 7\"],\n  [\"27\u00f74=6, 3\", \"10\u00f77=1, 3\"],\n  [\"99\u00f75=19, 4\", \"87\u00f74=21, 3\"],\n  [\"68\u00f75=13, 3\", \"37\u00f76=6, 1\"],\n  [\"80\u00f78=10, 0\", \"86\u00f76=14, 2\"],\n  [\"34\u00f76=5, 4\", \"66\u00f75=13, 1\"],\n  [\"42\u00f77=6, 0\", \"56\u00f78=7, 0\"],\n  [\"36\u00f76=6, 0\", \"34\u00f78=4, 2\"],\n  [\"10\u00f74=2, 2\", \"94\u00f77=13, 3\"],\n  [\"18\u00f72=9, 0\", \"17\u00f74=4, 1\"],\n  [\"70\u00f73=23, 1\", \"82\u00f78=10, 2\"],\n  [\"30\u00f72=15, 0\", \"57\u00f79=6, 3\"],\n  [\"79\u00f73=26, 1\", \"68\u00f78=8, 4\"],\n  [\"88\u00f79=9, 7\", \"45\u00f74=11, 1\"],\n  [\"11\u00f77=1, 4\", \"77\u00f75=15, 2\"],\n  [\"42\u00f78=5, 2\", \"38\u00f74=9, 2\"],\n  [\"97\u00f73=32, 1\", \"96\u00f75=19, 1\"],\n  [\"56\u00f73=18, 2\", \"92\u00f76=15, 2\"],\n  [\"78\u00f76=13, 0\", \"97\u00f74=24, 1\"],\n  [\"87\u00f76=14, 3\", \"67\u00f79=7, 4\"],\n  [\"18\u00f78=2, 2\", \"33\u00f72=16, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Collect the cell objects (row, col) in document order so that the\n// replacement at position N in `pairs` lines up with the Nth non-empty\n// cell encountered while scanning the table top-to-bottom, left-to-right\n// (duplicate old values, e.g. \"11\u00f77=1, 4\", are therefore resolved\n// positionally rather than by a naive global text replace).\nconst rows = table.values;\nlet pairIndex = 0;\nfor (let r = 0; r < rows.length; r++) {\n  const rowValues = rows[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    if (rowValues[c] === \"\") {\n      continue;\n    }\n    if (pairIndex >= pairs.length) {\n      break;\n    }\n    const [oldText, newText] = pairs[pairIndex];\n    if (rowValues[c] !== oldText) {\n      throw new Error(\n        \"Unexpected cell text at row \" + r + \", col \" + c +\n        \": got '\" + rowValues[c] + \"', expected '\" + oldText + \"'\"\n      );\n    }\n    const cell = table.getCell(r, c);\n    cell.value = newText;\n    pairIndex++;\n  }\n}\n\nawait context.sync();\n\nif (pairIndex !== pairs.length) {\n  throw new Error(\"Only replaced \" + pairIndex + \" of \" + pairs.length + \" cells\");\n}\n", "ps1": "# Update the answer text in each populated table cell of the practice\n# worksheet (\"two-digit number divided by one-digit number\") in document\n# order, mapping each old \"a\u00f7b=c, d\" string to its replacement.\n$pairs = @(\n    ,(\"11\u00f77=1, 4\", \"35\u00f77=5, 0\")\n    ,(\"14\u00f77=2, 0\", \"45\u00f73=15, 0\")\n    ,(\"89\u00f79=9, 8\", \"67\u00f79=7, 4\")\n    ,(\"60\u00f77=8, 4\", \"13\u00f79=1, 4\")\n    ,(\"66\u00f74=16, 2\", \"70\u00f79=7, 7\")\n    ,(\"27\u00f74=6, 3\", \"10\u00f77=1, 3\")\n    ,(\"99\u00f75=19, 4\", \"87\u00f74=21, 3\")\n    ,(\"68\u00f75=13, 3\", \"37\u00f76=6, 1\")\n    ,(\"80\u00f78=10, 0\", \"86\u00f76=14, 2\")\n    ,(\"34\u00f76=5, 4\", \"66\u00f75=13, 1\")\n    ,(\"42\u00f77=6, 0\", \"56\u00f78=7, 0\")\n    ,(\"36\u00f76=6, 0\", \"34\u00f78=4, 2\")\n    ,(\"10\u00f74=2, 2\", \"94\u00f77=13, 3\")\n    ,(\"18\u00f72=9, 0\", \"17\u00f74=4, 1\")\n    ,(\"70\u00f73=23, 1\", \"82\u00f78=10, 2\")\n    ,(\"30\u00f72=15, 0\", \"57\u00f79=6, 3\")\n    ,(\"79\u00f73=26, 1\", \"68\u00f78=8, 4\")\n    ,(\"88\u00f79=9, 7\", \"45\u00f74=11, 1\")\n    ,(\"11\u00f77=1, 4\", \"77\u00f75=15, 2\")\n    ,(\"42\u00f78=5, 2\", \"38\u00f74=9, 2\")\n    ,(\"97\u00f73=32, 1\", \"96\u00f75=19, 1\")\n    ,(\"56\u00f73=18, 2\", \"92\u00f76=15, 2\")\n    ,(\"78\u00f76=13, 0\", \"97\u00f74=24, 1\")\n    ,(\"87\u00f76=14, 3\", \"67\u00f79=7, 4\")\n    ,(\"18\u00f78=2, 2\", \"33\u00f72=16, 1\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$pairIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $rng = $cell.Range\n        # Cell range text includes trailing cell-mark control chars\n        # (CR + BEL); strip them to get the plain content.\n        $txt = $rng.Text -replace \"[`r`a]\", \"\"\n        if ($txt -eq \"\") {\n            continue\n        }\n        if ($pairIndex -ge $pairs.Length) {\n            break\n        }\n        $old = $pairs[$pairIndex][0]\n        $new = $pairs[$pairIndex][1]\n        if ($txt -ne $old) {\n            throw \"Unexpected cell text at row $r, col $c`: got '$txt', expected '$old'\"\n        }\n        $rng.Text = $new\n        $pairIndex++\n    }\n}\n\nif ($pairIndex -ne $pairs.Length) {\n    throw \"Only replaced $pairIndex of $($pairs.Length) cells\"\n}\n"}
